# Update column F (ticket/price-related numeric values) on the "展览" and
# "全部类型" worksheets to reflect the latest generated output.

$wb = $excel.ActiveWorkbook

# Map of row -> [old value, new value] is not needed; we just assign new
# values directly to column F for the affected rows on both sheets.

$updates = @{
    2  = 12857
    3  = 630
    5  = 37
    7  = 404
    9  = 12876
    10 = 40
    11 = 25
    12 = 5244
    13 = 546
    14 = 19
    15 = 13
    16 = 30
    18 = 38
    20 = 676
    21 = 2856
    22 = 6178
    23 = 1158
    24 = 3627
    26 = 45
}

$ws1 = $wb.Worksheets.Item("展览")
foreach ($row in $updates.Keys) {
    $ws1.Range("F$row").Value = $updates[$row]
}

$updates4 = @{
    2  = 12857
    3  = 630
    5  = 37
    8  = 404
    10 = 12876
    11 = 40
    12 = 25
    13 = 5244
    14 = 546
    15 = 19
    16 = 13
    17 = 30
    19 = 38
    21 = 676
    22 = 2856
    24 = 6178
    25 = 1158
    26 = 3627
    28 = 45
}

$ws4 = $wb.Worksheets.Item("全部类型")
foreach ($row in $updates4.Keys) {
    $ws4.Range("F$row").Value = $updates4[$row]
}
